$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserData")

# Update existing user data (row 2: Tom -> Tom3, email updated; row 3: john -> john4, email updated)
$ws.Range("A2").Value = "Tom3"
$ws.Range("D2").Value = "tomhanks3@zmail.com"
$ws.Range("A3").Value = "john4"
$ws.Range("D3").Value = "johnswam4@zmail.com"

# Widen column D to fit the longer email addresses
$ws.Columns.Item(4).ColumnWidth = 45.8

# Add the new (currently empty) "GetUserDetails" sheet after "UserData"
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "GetUserDetails"

# Keep "UserData" the active/visible sheet, with D9 as the selected cell
$ws.Activate()
$ws.Range("D9").Select() | Out-Null
